$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a text value to a cell while keeping it text (avoids
# Excel auto-converting numeric-looking strings like "302.38" into numbers)
# and without leaving the cells style index changed from the original.
function Set-TextValue($addr, $val) {
    $rng = $ws.Range($addr)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

# Row 2
Set-TextValue "D2" '42.941.62'
$ws.Range("E2").Value = '  -0.17%  '

# Row 3
Set-TextValue "D3" '2.333.18'
$ws.Range("E3").Value = '  +0.97%  '

# Row 4
$ws.Range("E4").Value = '  -0.01%  '

# Row 5
Set-TextValue "D5" '302.38'
$ws.Range("E5").Value = '  +0.13%  '

# Row 6
Set-TextValue "D6" '94.57'
$ws.Range("E6").Value = '  -2.60%  '

# Row 7
$ws.Range("E7").Value = '  -0.89%  '

# Row 8
$ws.Range("E8").Value = '  +0.00%  '

# Row 9
$ws.Range("E9").Value = '  -1.06%  '

# Row 10
Set-TextValue "D10" '33.95'
$ws.Range("E10").Value = '  -3.30%  '

# Row 11
Set-TextValue "D11" '0.0783'
$ws.Range("E11").Value = '  -1.09%  '

# Row 12
Set-TextValue "D12" '18.61'
$ws.Range("E12").Value = '  -4.01%  '

# Row 13
$ws.Range("E13").Value = '  +1.50%  '

# Row 14
Set-TextValue "D14" '6.73'
$ws.Range("E14").Value = '  -1.95%  '

# Row 15
Set-TextValue "D15" '2.699.58'
$ws.Range("E15").Value = '  +0.90%  '

# Row 16
Set-TextValue "D16" '2.311.62'
$ws.Range("E16").Value = '  +0.53%  '

# Row 17
$ws.Range("E17").Value = '  +0.86%  '

# Row 18
Set-TextValue "D18" '42.875.74'
$ws.Range("E18").Value = '  -0.19%  '

# Row 19
Set-TextValue "D19" '12.13'
$ws.Range("E19").Value = '  -3.11%  '

# Row 20
Set-TextValue "D20" '6.20'
$ws.Range("E20").Value = '  +2.63%  '

# Row 21
$ws.Range("E21").Value = '  -0.65%  '

# Row 22
Set-TextValue "D22" '67.83'
$ws.Range("E22").Value = '  +0.15%  '

# Row 23
Set-TextValue "D23" '235.25'
$ws.Range("E23").Value = '  -0.43%  '

# Row 24
$ws.Range("E24").Value = '  -0.89%  '

# Row 25
$ws.Range("E25").Value = '  -0.06%  '

# Row 26
$ws.Range("E26").Value = '  -1.54%  '

# Row 27
$ws.Range("E27").Value = '  -0.98%  '

# Row 28
Set-TextValue "D28" '2.36'
$ws.Range("E28").Value = '  +14.21%  '

# Row 29
Set-TextValue "D29" '9.17'
$ws.Range("E29").Value = '  +1.00%  '

# Row 30
Set-TextValue "D30" '31.42'
$ws.Range("E30").Value = '  -3.88%  '

# Row 31
$ws.Range("E31").Value = '  +0.01%  '

# Row 32
Set-TextValue "D32" '4.99'
$ws.Range("E32").Value = '  +0.09%  '

# Row 33
Set-TextValue "D33" '0.0728'
$ws.Range("E33").Value = '  +4.20%  '

# Row 34
Set-TextValue "D34" '17.18'
$ws.Range("E34").Value = '  -3.75%  '

# Row 35
$ws.Range("B35").Value = 'RenderToken'
$ws.Range("C35").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
Set-TextValue "D35" '4.37'
$ws.Range("E35").Value = '  -2.47%  '

# Row 36
$ws.Range("B36").Value = 'WEMIXToken'
$ws.Range("C36").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
Set-TextValue "D36" '2.32'
$ws.Range("E36").Value = '  -1.33%  '

# Row 37
$ws.Range("B37").Value = 'ARBITRUM'
$ws.Range("C37").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
Set-TextValue "D37" '1.82'
$ws.Range("E37").Value = '  +3.15%  '

# Row 38
$ws.Range("E38").Value = '  -0.20%  '

# Row 39
$ws.Range("B39").Value = 'LidoDAOToken'
$ws.Range("C39").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
Set-TextValue "D39" '2.74'
$ws.Range("E39").Value = '  -0.55%  '

# Row 40
$ws.Range("B40").Value = 'EnergySwap'
$ws.Range("C40").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
Set-TextValue "D40" '22.01'
$ws.Range("E40").Value = '  +16.49%  '

# Row 41
$ws.Range("B41").Value = 'Stellar'
$ws.Range("C41").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
Set-TextValue "D41" '0.108'
$ws.Range("E41").Value = '  -1.02%  '

# Row 42
$ws.Range("B42").Value = 'Monero'
$ws.Range("C42").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
Set-TextValue "D42" '111.91'
$ws.Range("E42").Value = '  -31.66%  '

# Row 43
Set-TextValue "D43" '1.930.65'
$ws.Range("E43").Value = '  -2.40%  '

# Row 44
Set-TextValue "D44" '0.0280'
$ws.Range("E44").Value = '  +0.23%  '

# Row 45
Set-TextValue "D45" '10.02'
$ws.Range("E45").Value = '  -5.60%  '

# Row 46
$ws.Range("E46").Value = '  +1.27%  '

# Row 47
Set-TextValue "D47" '2.71'
$ws.Range("E47").Value = '  -2.32%  '

# Row 48
Set-TextValue "D48" '2.564.62'
$ws.Range("E48").Value = '  +0.90%  '

# Row 49
Set-TextValue "D49" '2.82'
$ws.Range("E49").Value = '  -0.80%  '

# Row 50
Set-TextValue "D50" '52.87'
$ws.Range("E50").Value = '  -1.18%  '

# Row 51
Set-TextValue "D51" '71.74'
$ws.Range("E51").Value = '  -0.63%  '
